# Update the "K" (strikeouts) column (G) with regenerated values.
# This replaces the old "Strike#" derived values with the new K values
# per the commit message: "regen save_data to use K instead of Strike#,
# regen std/mean, calc and write s_vals"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 1
    5  = 0
    7  = 0
    8  = 0
    9  = 0
    10 = 1
    11 = 2
    12 = 1
    13 = 0
    14 = 2
    15 = 1
    16 = 2
    17 = 2
    18 = 1
    21 = 1
    22 = 1
    24 = 3
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
